$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles/number formats) from the last existing row (125)
# down into the two new rows so the appended rows match the existing look
# (bold/bordered/centered index column, date-formatted match-date column).
$ws.Range("A125:V125").Copy($ws.Range("A126:V126"))
$ws.Range("A125:V125").Copy($ws.Range("A127:V127"))

# Row 126
$ws.Range("A126").Value = 125
$ws.Range("B126").Value = "romania"
$ws.Range("C126").Value = "liga-1"
$ws.Range("D126").Value = "2023-2024"
$ws.Range("E126").Value = 45242.52083333334
$ws.Range("F126").Value = "FC Botosani"
$ws.Range("G126").Value = 3
$ws.Range("H126").Value = "FC Voluntari"
$ws.Range("I126").Value = 3
$ws.Range("J126").Value = 2.62
$ws.Range("K126").Value = "08/11/2023 08:12"
$ws.Range("L126").Value = 2.55
$ws.Range("M126").Value = "12/11/2023 12:22"
$ws.Range("N126").Value = 3.22
$ws.Range("O126").Value = "08/11/2023 08:12"
$ws.Range("P126").Value = 3.22
$ws.Range("Q126").Value = "12/11/2023 12:22"
$ws.Range("R126").Value = 2.78
$ws.Range("S126").Value = "08/11/2023 08:12"
$ws.Range("T126").Value = 2.94
$ws.Range("U126").Value = "12/11/2023 12:22"
$ws.Range("V126").Value = "https://www.betexplorer.com/football/romania/liga-1/fc-botosani-voluntari/2qw8eynN/"

# Row 127
$ws.Range("A127").Value = 126
$ws.Range("B127").Value = "romania"
$ws.Range("C127").Value = "liga-1"
$ws.Range("D127").Value = "2023-2024"
$ws.Range("E127").Value = 45242.625
$ws.Range("F127").Value = "Univ. Craiova"
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = "Din. Bucuresti"
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1.33
$ws.Range("K127").Value = "08/11/2023 08:12"
$ws.Range("L127").Value = 1.44
$ws.Range("M127").Value = "12/11/2023 14:56"
$ws.Range("N127").Value = 5.05
$ws.Range("O127").Value = "08/11/2023 08:12"
$ws.Range("P127").Value = 4.51
$ws.Range("Q127").Value = "12/11/2023 14:56"
$ws.Range("R127").Value = 9.24
$ws.Range("S127").Value = "08/11/2023 08:12"
$ws.Range("T127").Value = 7.49
$ws.Range("U127").Value = "12/11/2023 14:56"
$ws.Range("V127").Value = "https://www.betexplorer.com/football/romania/liga-1/univ-craiova-din-bucuresti/K0oiaJnb/"
